$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/17/2025  Through  11/23/2025"

# --- Cells that flip between numeric and placeholder-text ("0" / "***.*") ---
# Value is set first (numbers get a leading apostrophe to force text when needed),
# then number formatting is copied from a same-column donor cell that already has
# the desired style, so the resulting style index matches the rest of the column.
$ws.Range("C20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("F15").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("D15").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("E15").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remaining data-table numeric updates (rows 15-30) ---
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 7
$ws.Range("H15").Value = -71.428571428571
$ws.Range("I15").Value = 39
$ws.Range("J15").Value = 31
$ws.Range("K15").Value = 25.806451612903
$ws.Range("L15").Value = 30
$ws.Range("M15").Value = 105.263157894737
$ws.Range("N15").Value = -29.090909090909
# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 29
$ws.Range("H16").Value = -17.241379310344
$ws.Range("I16").Value = 413
$ws.Range("J16").Value = 426
$ws.Range("K16").Value = -3.051643192488
$ws.Range("L16").Value = 4.292929292929
$ws.Range("M16").Value = 14.088397790055
$ws.Range("N16").Value = -60.516252390057
# Row 17
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 41
$ws.Range("G17").Value = 53
$ws.Range("H17").Value = -22.641509433962
$ws.Range("I17").Value = 664
$ws.Range("J17").Value = 663
$ws.Range("K17").Value = 0.150829562594
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 93.586005830903
$ws.Range("N17").Value = 5.564387917329
# Row 18
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 30
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 76.470588235294
$ws.Range("I18").Value = 253
$ws.Range("J18").Value = 278
$ws.Range("K18").Value = -8.992805755395
$ws.Range("L18").Value = 13.452914798206
$ws.Range("M18").Value = 10.480349344978
$ws.Range("N18").Value = -78.631756756756
# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -8.333333333333
$ws.Range("I19").Value = 568
$ws.Range("J19").Value = 578
$ws.Range("K19").Value = -1.730103806228
$ws.Range("L19").Value = 18.333333333333
$ws.Range("M19").Value = 94.520547945205
$ws.Range("N19").Value = 35.885167464114
# Row 20
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = -53.571428571428
$ws.Range("J20").Value = 287
$ws.Range("K20").Value = -20.557491289198
$ws.Range("L20").Value = -22.972972972973
$ws.Range("M20").Value = 109.174311926605
$ws.Range("N20").Value = -58.695652173913
# Row 21
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = -24.444444444444
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 182
$ws.Range("H21").Value = -14.835164835164
$ws.Range("I21").Value = 2178
$ws.Range("J21").Value = 2273
$ws.Range("K21").Value = -4.179498460184
$ws.Range("L21").Value = 3.615604186489
$ws.Range("M21").Value = 60.382916053019
$ws.Range("N21").Value = -44.410413476263
# Row 22
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("M22").Value = 11.111111111111
# Row 23
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = -62.962962962963
# Row 24
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 49
$ws.Range("E24").Value = -51.020408163265
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 123
$ws.Range("H24").Value = -18.699186991869
$ws.Range("I24").Value = 1264
$ws.Range("J24").Value = 1026
$ws.Range("K24").Value = 23.196881091617
$ws.Range("L24").Value = 36.206896551724
$ws.Range("M24").Value = 44.292237442922
# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = -85.185185185185
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -40
$ws.Range("I25").Value = 327
$ws.Range("J25").Value = 324
$ws.Range("K25").Value = 0.925925925925
$ws.Range("L25").Value = 11.224489795918
# Row 26
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 93
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = 29.166666666666
$ws.Range("I26").Value = 969
$ws.Range("J26").Value = 915
$ws.Range("K26").Value = 5.901639344262
$ws.Range("L26").Value = 12.412993039443
$ws.Range("M26").Value = 9.244644870349
# Row 27
$ws.Range("D27").Value = 1
$ws.Range("I27").Value = 52
$ws.Range("J27").Value = 48
$ws.Range("K27").Value = 8.333333333333
$ws.Range("L27").Value = 1.960784313725
# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -16.666666666666
$ws.Range("I28").Value = 61
$ws.Range("J28").Value = 73
$ws.Range("K28").Value = -16.438356164383
$ws.Range("L28").Value = -10.294117647058
# Row 29
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("I29").Value = 33
$ws.Range("J29").Value = 39
$ws.Range("K29").Value = -15.384615384615
$ws.Range("L29").Value = -8.333333333333
$ws.Range("M29").Value = -15.384615384615
$ws.Range("N29").Value = -73.387096774193
# Row 30
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -75
$ws.Range("I30").Value = 27
$ws.Range("J30").Value = 34
$ws.Range("K30").Value = -20.588235294117
$ws.Range("L30").Value = 3.846153846153
$ws.Range("M30").Value = -15.625
$ws.Range("N30").Value = -74.038461538461
